$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "45"

$c9 = $ws.Range("C9")
$c9.Characters(48, 9).Text = "11/13/2022"
$c9.Characters(27, 10).Text = "11/7/2022"

# --- Cells that change from numeric to text (shared-string) ---
# Set apostrophe-prefixed text first (forces text type), then copy the
# number format from a neighboring cell that already has the General/text
# style (14) so the final style index matches the original text-style cells.
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "'***.*"
$ws.Range("F15").Value = "'0"
$ws.Range("D16").Value = "'0"
$ws.Range("E16").Value = "'***.*"
$ws.Range("D23").Value = "'0"
$ws.Range("E23").Value = "'***.*"
$ws.Range("F23").Value = "'0"
$ws.Range("D26").Value = "'0"
$ws.Range("E26").Value = "'***.*"
$ws.Range("F27").Value = "'0"

$ws.Range("C15").Copy() | Out-Null
$ws.Range("D15:F15").PasteSpecial(-4122) | Out-Null
$ws.Range("D16:E16").PasteSpecial(-4122) | Out-Null
$ws.Range("D23:F23").PasteSpecial(-4122) | Out-Null
$ws.Range("D26:E26").PasteSpecial(-4122) | Out-Null
$ws.Range("F27").PasteSpecial(-4122) | Out-Null

# --- Cells that change from text to numeric (need style/number format fix) ---
$ws.Range("C18").Value = 2
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("D20").Value = 1
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("E20").Value = 600
$ws.Range("E20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C26").Value = 1
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"

# --- Remaining pure numeric value updates (style unchanged) ---
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = -100
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 42.857142857142
$ws.Range("I16").Value = 68
$ws.Range("K16").Value = 19.298245614035
$ws.Range("L16").Value = 4.615384615384
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = -27.777777777777
$ws.Range("I17").Value = 189
$ws.Range("J17").Value = 146
$ws.Range("K17").Value = 29.452054794520
$ws.Range("L17").Value = 13.173652694610
$ws.Range("E18").Value = 100
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = -20
$ws.Range("I18").Value = 75
$ws.Range("J18").Value = 73
$ws.Range("K18").Value = 2.739726027397
$ws.Range("L18").Value = 2.739726027397
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = 45.833333333333
$ws.Range("I19").Value = 352
$ws.Range("J19").Value = 274
$ws.Range("K19").Value = 28.467153284671
$ws.Range("L19").Value = 44.855967078189
$ws.Range("C20").Value = 7
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 500
$ws.Range("I20").Value = 98
$ws.Range("J20").Value = 44
$ws.Range("K20").Value = 122.727272727273
$ws.Range("L20").Value = 113.04347826087
$ws.Range("C21").Value = 19
$ws.Range("E21").Value = 18.75
$ws.Range("F21").Value = 80
$ws.Range("G21").Value = 58
$ws.Range("H21").Value = 37.931034482758
$ws.Range("I21").Value = 792
$ws.Range("J21").Value = 613
$ws.Range("K21").Value = 29.200652528548
$ws.Range("L21").Value = 28.571428571428
$ws.Range("H23").Value = -100
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = 20
$ws.Range("F24").Value = 95
$ws.Range("G24").Value = 76
$ws.Range("H24").Value = 25
$ws.Range("I24").Value = 1219
$ws.Range("J24").Value = 801
$ws.Range("K24").Value = 52.184769038701
$ws.Range("L24").Value = 67.215363511659
$ws.Range("C25").Value = 11
$ws.Range("E25").Value = 22.222222222222
$ws.Range("F25").Value = 40
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = 11.111111111111
$ws.Range("I25").Value = 435
$ws.Range("J25").Value = 355
$ws.Range("K25").Value = 22.535211267605
$ws.Range("L25").Value = 40.776699029126
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -75
$ws.Range("I26").Value = 17
$ws.Range("K26").Value = -39.285714285714
$ws.Range("L26").Value = -19.047619047619
$ws.Range("D27").Value = 1
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -100
$ws.Range("J27").Value = 42
$ws.Range("K27").Value = -9.523809523809
$ws.Range("J30").Value = 3
$ws.Range("K30").Value = 166.666666666667
